$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows immediately above the current row 914, shifting the
# existing rows 914-1013 down to become rows 917-1016 (data unchanged).
$ws.Rows("914:916").Insert()

# ---- New row 914: Especial / Peru batch -------------------------------
$ws.Range("A914").Value = 4
$ws.Range("B914").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C914").Value = "Los Lagos"
$ws.Range("D914").Value = 45142
$ws.Range("E914").Value = 10
$ws.Range("F914").Value = "Fruta"
$ws.Range("G914").Value = 100106
$ws.Range("H914").Value = "Oleaginosos"
$ws.Range("I914").Value = 100106002
$ws.Range("J914").Value = "Palta"
$ws.Range("K914").Value = "Hass"
$ws.Range("L914").Value = "Especial"
$ws.Range("M914").Value = 150
$ws.Range("N914").Value = 38000
$ws.Range("O914").Value = 38000
$ws.Range("P914").Value = 38000
$ws.Range("Q914").Value = "`$/bandeja 10 kilos"
$ws.Range("R914").Value = "Perú"
$ws.Range("S914").Value = 3800
$ws.Range("T914").Value = 10

# ---- New row 915: Primera / Peru batch --------------------------------
$ws.Range("A915").Value = 4
$ws.Range("B915").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C915").Value = "Los Lagos"
$ws.Range("D915").Value = 45142
$ws.Range("E915").Value = 10
$ws.Range("F915").Value = "Fruta"
$ws.Range("G915").Value = 100106
$ws.Range("H915").Value = "Oleaginosos"
$ws.Range("I915").Value = 100106002
$ws.Range("J915").Value = "Palta"
$ws.Range("K915").Value = "Hass"
$ws.Range("L915").Value = "Primera"
$ws.Range("M915").Value = 200
$ws.Range("N915").Value = 32000
$ws.Range("O915").Value = 32000
$ws.Range("P915").Value = 32000
$ws.Range("Q915").Value = "`$/bandeja 10 kilos"
$ws.Range("R915").Value = "Perú"
$ws.Range("S915").Value = 3200
$ws.Range("T915").Value = 10

# ---- New row 916: Segunda / Peru batch --------------------------------
$ws.Range("A916").Value = 4
$ws.Range("B916").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C916").Value = "Los Lagos"
$ws.Range("D916").Value = 45142
$ws.Range("E916").Value = 10
$ws.Range("F916").Value = "Fruta"
$ws.Range("G916").Value = 100106
$ws.Range("H916").Value = "Oleaginosos"
$ws.Range("I916").Value = 100106002
$ws.Range("J916").Value = "Palta"
$ws.Range("K916").Value = "Hass"
$ws.Range("L916").Value = "Segunda"
$ws.Range("M916").Value = 200
$ws.Range("N916").Value = 28000
$ws.Range("O916").Value = 28000
$ws.Range("P916").Value = 28000
$ws.Range("Q916").Value = "`$/bandeja 10 kilos"
$ws.Range("R916").Value = "Perú"
$ws.Range("S916").Value = 2800
$ws.Range("T916").Value = 10
